# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.692.33'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '3.144.76'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''534.75'
$ws.Range('E5').Value = '  +1.33%  '
$ws.Range('D6').Value = '''143.47'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.142.63'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D10').Value = '''7.18'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '''0.396'
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('D13').Value = '3.682.67'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').Value = '''25.81'
$ws.Range('E15').Value = '  -4.07%  '
$ws.Range('D16').Value = '''0.0000167'
$ws.Range('D17').Value = '58.694.74'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').Value = '3.135.28'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '''6.14'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = '''12.94'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '''8.01'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').Value = '''343.86'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '''0.515'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').Value = '''67.95'
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '0.0₃0938'
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('D29').Value = '''7.57'
$ws.Range('E29').Value = '  +4.10%  '
$ws.Range('D30').Value = '''6.50'
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '''1.91'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '''21.23'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '''1.21'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '''4.82'
$ws.Range('E35').Value = '  +3.25%  '
$ws.Range('D36').Value = '''158.11'
$ws.Range('E36').Value = '  +2.53%  '
$ws.Range('D37').Value = '''6.29'
$ws.Range('E37').Value = '  +3.51%  '
$ws.Range('D38').Value = '''26.29'
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('D39').Value = '''1.27'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').Value = '''1.67'
$ws.Range('E40').Value = '  +12.03%  '
$ws.Range('D41').Value = '''0.0675'
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('D43').Value = '''4.03'
$ws.Range('D44').Value = '3.181.32'
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').Value = '''36.62'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = '''0.999'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').Value = '''0.0266'
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('D48').Value = '2.311.17'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('E49').Value = '  +4.91%  '
$ws.Range('D50').Value = '''20.78'
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').Value = '''6.11'
$ws.Range('E51').Value = '  +1.94%  '
